$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column D (shifts existing D:K to E:L) for the new 2018-12-31 year
$ws.Columns.Item(4).EntireColumn.Insert()

# Copy number formats/styles from the (old) data column now in E into the new D column
# for each of the three contiguous data blocks on the sheet
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the 2018-12-31 figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 5506200
$ws.Range("D9").Value = 4527900
$ws.Range("D10").Value = 978300
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 70500
$ws.Range("D17").Value = 5303300
$ws.Range("D18").Value = 202900
$ws.Range("D20").Value = 1400
$ws.Range("D21").Value = 389300
$ws.Range("D22").Value = 21100
$ws.Range("D23").Value = 183200
$ws.Range("D24").Value = 44100
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 139100
$ws.Range("D27").Value = 139100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -1400
$ws.Range("D33").Value = 139100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 139100

$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 131700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 203500
$ws.Range("D44").Value = 1339900
$ws.Range("D45").Value = 12800
$ws.Range("D46").Value = 1687900
$ws.Range("D47").Value = "NA"
$ws.Range("D48").Value = 1184100
$ws.Range("D49").Value = 291400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 38000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3201300
$ws.Range("D57").Value = 127500
$ws.Range("D58").Value = 1204600
$ws.Range("D59").Value = 161200
$ws.Range("D60").Value = 1493300
$ws.Range("D61").Value = 488700
$ws.Range("D62").Value = 152400
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 2134400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 942300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1066900
$ws.Range("D77").Value = 0

$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 139100
$ws.Range("D83").Value = 185100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 215400
$ws.Range("D91").Value = -238300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -227200
$ws.Range("D96").Value = -9300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 19100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 7200
